$wb = $excel.ActiveWorkbook

# ---- Sheet "Tiến độ" ----
$ws = $wb.Worksheets.Item("Tiến độ")

$workItems = "Hoàn thành giao diện đăng nhập`nHoàn thành giao diện trang sản phẩm`nHoàn thành giao diện danh mục sản phẩm cần mua`nTìm được dữ liệu mẫu cho ứng dụng`nHoàn thành chức năng đăng nhập`nHoàn thành chức năng thêm sản phẩm vào danh mục sản phẩm cần mua`nHoàn thành chức năng xác nhận danh mục sản phẩm cần mua`nHoàn thành chức năng xem sản phẩm"

# Clear old per-row notes in column C (rows 3-4) that are being consolidated
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""

# Row 3 content
$ws.Range("A3").Value = "30/1/2018 - 9/2/2018"
$ws.Range("B3").Value = $workItems
$ws.Range("D3").Value = $workItems
$ws.Range("E3").Value = 0.8
$ws.Range("F3").Value = "Trễ 5 ngày so với tiến độ"

# Merge B3:C3 and set alignment/wrap
$ws.Range("B3:C3").Merge()
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").WrapText = $true
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("D3").WrapText = $true
$ws.Range("E3").NumberFormat = "0%"

# Box borders on A3,D3,E3,F3; split box on merged B3:C3
foreach ($addr in @("A3","D3","E3","F3")) {
  $c = $ws.Range($addr)
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(10).LineStyle = 1
}
$b3 = $ws.Range("B3")
$b3.Borders.Item(7).LineStyle = 1
$b3.Borders.Item(8).LineStyle = 1
$b3.Borders.Item(9).LineStyle = 1
$c3 = $ws.Range("C3")
$c3.Borders.Item(8).LineStyle = 1
$c3.Borders.Item(9).LineStyle = 1
$c3.Borders.Item(10).LineStyle = 1

# Row height for the now-tall row 3
$ws.Rows.Item(3).RowHeight = 148.5

# Column widths
$ws.Columns.Item(1).ColumnWidth = 18.736979166666668
$ws.Columns.Item(3).ColumnWidth = 57.307291666666664
$ws.Columns.Item(4).ColumnWidth = 54.877604166666664
$ws.Columns.Item(6).ColumnWidth = 21.877604166666668

# Page setup + selection
$ws.PageSetup.Orientation = 1
$ws.Activate()
$ws.Range("F9").Select()

# ---- Sheet "Danh sách công việc" ----
$ws1 = $wb.Worksheets.Item("Danh sách công việc")
$ws1.Activate()
$ws1.Range("B12").Select()

# Re-activate "Tiến độ" as the selected tab (matches tabSelected in source)
$ws.Activate()
